$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '24.687.44'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.64%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.696.20'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.74%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9993'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.27%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '316.68'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.03%  '
$ws.Range("E6").Value = '  +0.17%  '
$ws.Range("E7").Value = '  +1.59%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4026'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.78%  '
$ws.Range("E9").Value = '  +2.85%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9997'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.31%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '51.79'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.24%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08763'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.96%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.230'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +5.43%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '23.43'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.91%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.113'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +10.64%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001317'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.52%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.692.78'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.62%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '99.75'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.48%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.07040'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.17%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '19.78'
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.085'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +7.01%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9995'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.26%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '14.30'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.44%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '24.682.51'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.65%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.155'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +7.74%  '
$ws.Range("E26").Value = '  +1.41%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.92'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +5.18%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '162.44'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.39%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '137.27'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +5.20%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.225'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.84%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.497'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +4.10%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.884.83'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.87%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.080'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.07%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08655'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.90%  '
$ws.Range("B35").Value = 'InternetComputer(DFINITY)'
$ws.Range("C35").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '7.137'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +6.08%  '
$ws.Range("B36").Value = 'FraxShare'
$ws.Range("C36").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '11.59'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +10.38%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.2753'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.62%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.927'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.54%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '14.50'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.53%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.09145'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.36%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.02726'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +7.29%  '
$ws.Range("E42").Value = '  +2.09%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.7668'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.23%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.642'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +9.15%  '
$ws.Range("E45").Value = '  +1.33%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '15.68'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.27%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.228'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.98%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.9989'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.19%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '141.03'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.95%  '
$ws.Range("E50").Value = '  +9.76%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07997'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.01%  '
